$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.277.37'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -1.89%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.904.94'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -2.66%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9991'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.10%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '333.83'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -2.38%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4635'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -3.05%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.4140'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.24%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '47.91'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.01%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.08040'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -2.72%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.020'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -2.00%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '22.21'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -2.73%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.962.16'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.80%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.960'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -3.62%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.133'
$ws.Range('D15').Style = 'Normal'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '89.20'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -3.18%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.9996'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.12%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.00001032'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -2.69%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06588'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.67%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '17.67'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -2.25%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.9992'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.07%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '29.197.20'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -2.02%  '
$ws.Range('E23').Value = '  -1.65%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.42'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.12%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.196'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -3.68%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.084.17'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -4.35%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '157.19'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -2.50%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '19.83'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -2.29%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.143'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -2.01%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '5.675'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.58%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '117.29'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -4.88%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.043'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +3.41%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.09454'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -2.04%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.428'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -3.64%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.549'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -3.67%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.385'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -2.49%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.06111'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -2.75%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.02260'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -2.44%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '8.442'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.36%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.180'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.41%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.5882'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -3.66%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.9990'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.11%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.1831'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -3.55%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '10.19'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -5.47%  '
$ws.Range('B45').Value = 'RenderToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.354'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -2.01%  '
$ws.Range('B46').Value = 'WEMIXTOKEN'
$ws.Range('C46').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.235'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -3.48%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.07520'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +1.78%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.5567'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -2.83%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '12.18'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -2.56%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.931'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -3.21%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '113.22'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.26%  '
